# Reformat works-cited URLs as clean hyperlinked text instead of raw URLs.
# For each citation paragraph, the trailing run containing the raw GitHub
# URL is converted into a real w:hyperlink (with an External relationship)
# whose visible text is a short label ("[View Report]", "[View 10-K]" or
# "[View Data]"), while keeping the original character formatting
# (Times New Roman, size 24 half-points, blue color, single underline).

$d = $word.ActiveDocument

$links = @(
    @{ Url = "https://github.com/kh0pper/DSCI-5330-Assignment-02/blob/main/Annual%20Report/2016-annual-report.pdf"; Text = "[View Report]" },
    @{ Url = "https://github.com/kh0pper/DSCI-5330-Assignment-02/blob/main/Annual%20Report/Final-Annual-Report-2017.pdf"; Text = "[View Report]" },
    @{ Url = "https://github.com/kh0pper/DSCI-5330-Assignment-02/blob/main/Annual%20Report/2018-Annual-Report.pdf"; Text = "[View Report]" },
    @{ Url = "https://github.com/kh0pper/DSCI-5330-Assignment-02/blob/main/Annual%20Report/Ford-2019-Printed-Annual-Report.pdf"; Text = "[View Report]" },
    @{ Url = "https://github.com/kh0pper/DSCI-5330-Assignment-02/blob/main/Annual%20Report/Ford-2020-Annual-Report-April-2020.pdf"; Text = "[View Report]" },
    @{ Url = "https://github.com/kh0pper/DSCI-5330-Assignment-02/blob/main/Annual%20Report/Ford-2021-Annual-Report.pdf"; Text = "[View Report]" },
    @{ Url = "https://github.com/kh0pper/DSCI-5330-Assignment-02/blob/main/Annual%20Report/2022-Annual-Report-1.pdf"; Text = "[View Report]" },
    @{ Url = "https://github.com/kh0pper/DSCI-5330-Assignment-02/blob/main/Annual%20Report/2023-Ford-Annual-Report.pdf"; Text = "[View Report]" },
    @{ Url = "https://github.com/kh0pper/DSCI-5330-Assignment-02/blob/main/Annual%20Report/Ford-2024-Annual-Report.pdf"; Text = "[View Report]" },
    @{ Url = "https://github.com/kh0pper/DSCI-5330-Assignment-02/blob/main/10k/2016_10K_for%20Year%20End%202015%20-%20filed%2002.11.16.pdf"; Text = "[View 10-K]" },
    @{ Url = "https://github.com/kh0pper/DSCI-5330-Assignment-02/blob/main/10k/2017_10K_for%20Year%20End%202016%20-%20filed%2002.09.17.pdf"; Text = "[View 10-K]" },
    @{ Url = "https://github.com/kh0pper/DSCI-5330-Assignment-02/blob/main/10k/2018_10K_for%20Year%20End%202017%20-%20filed%2002.08.18.pdf"; Text = "[View 10-K]" },
    @{ Url = "https://github.com/kh0pper/DSCI-5330-Assignment-02/blob/main/10k/2019_10K_for%20Year%20End%202018%20-%20filed%2002.21.19.pdf"; Text = "[View 10-K]" },
    @{ Url = "https://github.com/kh0pper/DSCI-5330-Assignment-02/blob/main/10k/2020_10K_for%20Year%20End%202019%20-%20filed%2002.05.20.pdf"; Text = "[View 10-K]" },
    @{ Url = "https://github.com/kh0pper/DSCI-5330-Assignment-02/blob/main/10k/2021_10K_for%20Year%20End%202020%20-%20filed%2002.05.21.pdf"; Text = "[View 10-K]" },
    @{ Url = "https://github.com/kh0pper/DSCI-5330-Assignment-02/blob/main/10k/2022_10K_for%20Year%20End%202021%20-%20filed%2002.04.22.pdf"; Text = "[View 10-K]" },
    @{ Url = "https://github.com/kh0pper/DSCI-5330-Assignment-02/blob/main/10k/2023_10K_for%20Year%20End%202022%20-%20filed%2002.03.23.pdf"; Text = "[View 10-K]" },
    @{ Url = "https://github.com/kh0pper/DSCI-5330-Assignment-02/blob/main/10k/2024_10K_for%20Year%20End%202023%20-%20filed%2002.07.24.pdf"; Text = "[View 10-K]" },
    @{ Url = "https://github.com/kh0pper/DSCI-5330-Assignment-02/blob/main/10k/2025_10K_for%20Year%20End%202024%20-%20filed%2002.06.25.pdf"; Text = "[View 10-K]" },
    @{ Url = "https://github.com/kh0pper/DSCI-5330-Assignment-02/blob/main/Ford_10K_Financial_Ratios_2015_2024.xlsx"; Text = "[View Data]" }
)

foreach ($link in $links) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $found = $rng.Find.Execute($link.Url, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "NOT FOUND: $($link.Url)"
        continue
    }

    # Replace the URL text with the short label before converting to a
    # hyperlink, so Hyperlinks.Add stores the label as TextToDisplay.
    $rng.Text = $link.Text

    $h = $d.Hyperlinks.Add($rng, $link.Url, "", "", $link.Text)

    # Re-apply the original direct character formatting (Hyperlinks.Add
    # replaces it with the built-in "Hyperlink" character style only).
    $hr = $h.Range
    $hr.Font.Name = "Times New Roman"
    $hr.Font.Size = 12
    $hr.Font.Color = 16711680
    $hr.Font.Underline = 1
}

Write-Host "Done converting $($links.Count) citation links to hyperlinks."
